$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed date) column C for rows 2-6 from 45212 to 45221
$ws.Range("C2:C6").Value = 45221
